# Auto-generated script to apply numeric updates to Famfrit_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 5278.4287
$ws.Range("I125").Value = 4200
$ws.Range("J125").Value = 5458.1665
$ws.Range("K125").Value = 37800
$ws.Range("L125").Value = 49123.4985
$ws.Range("M125").Value = -35340
$ws.Range("N125").Value = -54043.4985
$ws.Range("H132").Value = 3213.3845
$ws.Range("I132").Value = 3564.45
$ws.Range("J132").Value = 2043.1666
$ws.Range("K132").Value = 10693.35
$ws.Range("L132").Value = 6129.4998
$ws.Range("M132").Value = -8163.349999999999
$ws.Range("N132").Value = -11189.4998
$ws.Range("H137").Value = 2899.75
$ws.Range("J137").Value = 1857.375
$ws.Range("L137").Value = 5572.125
$ws.Range("N137").Value = -10672.125
$ws.Range("H141").Value = 2744.4
$ws.Range("I141").Value = 2797.5715
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 8392.7145
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = -3212.7145
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10874898
$ws.Range("I32").Value = 13515532
$ws.Range("K32").Value = 13515532
$ws.Range("M32").Value = -13515245
$ws.Range("H63").Value = 4147.486
$ws.Range("I63").Value = 2692.1482
$ws.Range("K63").Value = 2692.1482
$ws.Range("M63").Value = -2006.1482
$ws.Range("H66").Value = 4147.486
$ws.Range("I66").Value = 2692.1482
$ws.Range("K66").Value = 13460.741
$ws.Range("M66").Value = -10028.741
$ws.Range("H88").Value = 7543.6665
$ws.Range("I88").Value = 13080.223
$ws.Range("J88").Value = 2007.1111
$ws.Range("K88").Value = 13080.223
$ws.Range("L88").Value = 2007.1111
$ws.Range("M88").Value = -12674.223
$ws.Range("N88").Value = -2819.1111
$ws.Range("H91").Value = 7543.6665
$ws.Range("I91").Value = 13080.223
$ws.Range("J91").Value = 2007.1111
$ws.Range("K91").Value = 13080.223
$ws.Range("L91").Value = 2007.1111
$ws.Range("M91").Value = -11676.223
$ws.Range("N91").Value = -4815.1111
$ws.Range("H129").Value = 44750
$ws.Range("J129").Value = 44750
$ws.Range("L129").Value = 44750
$ws.Range("N129").Value = -54750
$ws.Range("H132").Value = 23811464
$ws.Range("I132").Value = 1756.8975
$ws.Range("J132").Value = 333337660
$ws.Range("K132").Value = 5270.6925
$ws.Range("L132").Value = 1000012980
$ws.Range("M132").Value = -2740.6925
$ws.Range("N132").Value = -1000018040

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2115.3
$ws.Range("I20").Value = 2054.4666
$ws.Range("K20").Value = 2054.4666
$ws.Range("M20").Value = -1807.4666
$ws.Range("H64").Value = 1834.4706
$ws.Range("I64").Value = 1588.2
$ws.Range("K64").Value = 1588.2
$ws.Range("M64").Value = -1363.2
$ws.Range("H67").Value = 1834.4706
$ws.Range("I67").Value = 1588.2
$ws.Range("K67").Value = 1588.2
$ws.Range("M67").Value = -808.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5269250
$ws.Range("I4").Value = 7010666.5
$ws.Range("J4").Value = 45000
$ws.Range("K4").Value = 7010666.5
$ws.Range("L4").Value = 45000
$ws.Range("M4").Value = -7010554.5
$ws.Range("N4").Value = -45224
$ws.Range("H31").Value = 18522728
$ws.Range("I31").Value = 3034.3684
$ws.Range("J31").Value = 62507000
$ws.Range("K31").Value = 3034.3684
$ws.Range("L31").Value = 62507000
$ws.Range("M31").Value = -2739.3684
$ws.Range("N31").Value = -62507590
$ws.Range("H34").Value = 18522728
$ws.Range("I34").Value = 3034.3684
$ws.Range("J34").Value = 62507000
$ws.Range("K34").Value = 3034.3684
$ws.Range("L34").Value = 62507000
$ws.Range("M34").Value = -2832.3684
$ws.Range("N34").Value = -62507404
$ws.Range("H58").Value = 2890.5
$ws.Range("I58").Value = 2543.9285
$ws.Range("J58").Value = 4103.5
$ws.Range("K58").Value = 2543.9285
$ws.Range("L58").Value = 4103.5
$ws.Range("M58").Value = -2340.9285
$ws.Range("N58").Value = -4509.5
$ws.Range("H86").Value = 4223.8
$ws.Range("I86").Value = 4029.875
$ws.Range("K86").Value = 4029.875
$ws.Range("M86").Value = -2906.875
$ws.Range("H89").Value = 4223.8
$ws.Range("I89").Value = 4029.875
$ws.Range("K89").Value = 20149.375
$ws.Range("M89").Value = -14533.375
$ws.Range("H99").Value = 15922.913
$ws.Range("I99").Value = 19556.6
$ws.Range("J99").Value = 14913.556
$ws.Range("K99").Value = 19556.6
$ws.Range("L99").Value = 14913.556
$ws.Range("M99").Value = -18058.6
$ws.Range("N99").Value = -17909.556
$ws.Range("H126").Value = 15922.913
$ws.Range("I126").Value = 19556.6
$ws.Range("J126").Value = 14913.556
$ws.Range("K126").Value = 58669.8
$ws.Range("L126").Value = 44740.66800000001
$ws.Range("M126").Value = -56199.8
$ws.Range("N126").Value = -49680.66800000001
$ws.Range("H132").Value = 3720.9048
$ws.Range("I132").Value = 2976.889
$ws.Range("K132").Value = 8930.667000000001
$ws.Range("M132").Value = -6400.667000000001
$ws.Range("H136").Value = 2890.5
$ws.Range("I136").Value = 2543.9285
$ws.Range("J136").Value = 4103.5
$ws.Range("K136").Value = 7631.7855
$ws.Range("L136").Value = 12310.5
$ws.Range("M136").Value = -5081.7855
$ws.Range("N136").Value = -17410.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4700
$ws.Range("I3").Value = 4700
$ws.Range("K3").Value = 14100
$ws.Range("M3").Value = -13988
$ws.Range("H34").Value = 2053.2856
$ws.Range("J34").Value = 2855.3333
$ws.Range("L34").Value = 8565.999899999999
$ws.Range("N34").Value = -8733.999899999999
$ws.Range("H38").Value = 98.23529000000001
$ws.Range("I38").Value = 78
$ws.Range("J38").Value = 127.14286
$ws.Range("K38").Value = 234
$ws.Range("L38").Value = 381.42858
$ws.Range("M38").Value = 113
$ws.Range("N38").Value = -1075.42858
$ws.Range("H39").Value = 3223.1
$ws.Range("J39").Value = 4171.2856
$ws.Range("L39").Value = 12513.8568
$ws.Range("N39").Value = -13101.8568
$ws.Range("H55").Value = 4000862.2
$ws.Range("J55").Value = 962.9048
$ws.Range("L55").Value = 2888.7144
$ws.Range("N55").Value = -3242.7144
$ws.Range("H68").Value = 1500
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 6000
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -7622
$ws.Range("H71").Value = 1500
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 18000
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -26112
$ws.Range("H118").Value = 898.9091
$ws.Range("I118").Value = 898.9091
$ws.Range("K118").Value = 2696.7273
$ws.Range("M118").Value = -1453.7273
$ws.Range("H122").Value = 1633.7693
$ws.Range("J122").Value = 1671.2222
$ws.Range("L122").Value = 15040.9998
$ws.Range("N122").Value = -19940.9998
$ws.Range("H129").Value = 3905.4375
$ws.Range("J129").Value = 3723.9583
$ws.Range("L129").Value = 11171.8749
$ws.Range("N129").Value = -21171.8749
$ws.Range("H132").Value = 4765093.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4765093.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 42885841.5
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -42890901.5
$ws.Range("H141").Value = 11657.889
$ws.Range("I141").Value = 10073.637
$ws.Range("K141").Value = 30220.911
$ws.Range("M141").Value = -25040.911

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1501.8235
$ws.Range("J97").Value = 3017.2856
$ws.Range("L97").Value = 3017.2856
$ws.Range("N97").Value = -4009.2856
$ws.Range("H126").Value = 20007266
$ws.Range("I126").Value = 16676500
$ws.Range("J126").Value = 22227778
$ws.Range("K126").Value = 50029500
$ws.Range("L126").Value = 66683334
$ws.Range("M126").Value = -50027030
$ws.Range("N126").Value = -66688274
$ws.Range("H132").Value = 2996.077
$ws.Range("I132").Value = 2954.2917
$ws.Range("K132").Value = 8862.875100000001
$ws.Range("M132").Value = -6332.875100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3775.4736
$ws.Range("I40").Value = 2949.818
$ws.Range("K40").Value = 2949.818
$ws.Range("M40").Value = -2813.818
$ws.Range("H136").Value = 9995
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 202001100
$ws.Range("I100").Value = 336667170
$ws.Range("J100").Value = 1999.5
$ws.Range("K100").Value = 673334340
$ws.Range("L100").Value = 3999
$ws.Range("M100").Value = -673333799
$ws.Range("N100").Value = -5081
$ws.Range("H132").Value = 5495
$ws.Range("I132").Value = 5313.636
$ws.Range("K132").Value = 15940.908
$ws.Range("M132").Value = -13410.908

Write-Output "Applied: 221 set, 6 added, 2 deleted"